$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.988.55'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.957.40'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.50%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.11'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.88%  '
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4869'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.78%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2947'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07038'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +3.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.77'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '107.58'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.963.31'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07822'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.473'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7013'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '281.07'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -4.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '31.002.03'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.34'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.61%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007808'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.222.36'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.572'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.91%  '
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.508'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.840'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.93'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.94'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.89%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.187'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1049'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.385'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -4.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.624'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -4.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.572'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.456'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04911'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7521'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.169'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.736'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02008'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.684'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.543'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.33%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '77.88'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +9.16%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.129'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9019'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.65%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '109.24'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.96%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4450'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.101'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +7.77%  '
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '995.78'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +7.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.380'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1251'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.87'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.46%  '
